$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.780.18"
$ws.Range("D3").Value = "2.074.96"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'232.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'58.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.0784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "'14.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.381.99"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "2.079.50"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "37.711.66"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "'6.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("D20").Value = "'71.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").Value = "'229.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'9.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.34%  "
$ws.Range("D27").Value = "'171.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'19.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "'4.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("D33").Value = "'0.0631"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "'2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D37").Value = "'3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'5.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("E40").Value = "  +7.37%  "
$ws.Range("D41").Value = "'101.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.0973"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'2.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "'17.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.53%  "
$ws.Range("D45").Value = "1.452.04"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("D48").Value = "'4.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.24%  "
$ws.Range("D49").Value = "'7.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "2.268.29"
$ws.Range("E51").Value = "  -0.49%  "
